$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.31"
$ws.Range("E2").Value = "'-3.35%"
$ws.Range("G2").Value = "'19"
$ws.Range("D3").Value = "'41.82"
$ws.Range("E3").Value = "'-5.81%"
$ws.Range("G3").Value = "'19"
$ws.Range("D4").Value = "'5.187"
$ws.Range("E4").Value = "'0.55%"
$ws.Range("G4").Value = "'19"
$ws.Range("D5").Value = "'0.08045"
$ws.Range("E5").Value = "'-4.01%"
$ws.Range("G5").Value = "'19"
$ws.Range("D6").Value = "'4.369"
$ws.Range("E6").Value = "'-1.73%"
$ws.Range("G6").Value = "'19"
$ws.Range("D7").Value = "'1.741"
$ws.Range("E7").Value = "'-10.91%"
$ws.Range("G7").Value = "'19"
$ws.Range("D8").Value = "'0.9277"
$ws.Range("E8").Value = "'-4.80%"
$ws.Range("G8").Value = "'19"
$ws.Range("D9").Value = "'0.1132"
$ws.Range("E9").Value = "'0.26%"
$ws.Range("G9").Value = "'19"
$ws.Range("D10").Value = "'0.1851"
$ws.Range("E10").Value = "'-2.36%"
$ws.Range("G10").Value = "'19"
$ws.Range("D11").Value = "'0.09333"
$ws.Range("E11").Value = "'-3.44%"
$ws.Range("G11").Value = "'19"
$ws.Range("E12").Value = "'-1.25%"
$ws.Range("G12").Value = "'19"
$ws.Range("D13").Value = "'7.367"
$ws.Range("E13").Value = "'-15.25%"
$ws.Range("G13").Value = "'19"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("E14").Value = "'-0.60%"
$ws.Range("G14").Value = "'19"
$ws.Range("D15").Value = "'0.001278"
$ws.Range("E15").Value = "'-0.92%"
$ws.Range("G15").Value = "'19"
$ws.Range("D16").Value = "'0.005950"
$ws.Range("E16").Value = "'3.81%"
$ws.Range("G16").Value = "'19"
$ws.Range("E17").Value = "'-1.32%"
$ws.Range("G17").Value = "'19"
$ws.Range("D18").Value = "'2.570"
$ws.Range("E18").Value = "'1.50%"
$ws.Range("G18").Value = "'19"
$ws.Range("D19").Value = "'0.3385"
$ws.Range("E19").Value = "'0.65%"
$ws.Range("G19").Value = "'19"
$ws.Range("E20").Value = "'1.45%"
$ws.Range("G20").Value = "'19"
$ws.Range("E21").Value = "'0.96%"
$ws.Range("G21").Value = "'19"
$ws.Range("D22").Value = "'0.04172"
$ws.Range("E22").Value = "'0.39%"
$ws.Range("G22").Value = "'19"
$ws.Range("D23").Value = "'0.001249"
$ws.Range("E23").Value = "'1.14%"
$ws.Range("G23").Value = "'19"
$ws.Range("D24").Value = "'0.004277"
$ws.Range("E24").Value = "'-3.03%"
$ws.Range("G24").Value = "'19"
$ws.Range("D25").Value = "'0.0001224"
$ws.Range("E25").Value = "'-6.02%"
$ws.Range("G25").Value = "'19"
$ws.Range("D26").Value = "'0.0002988"
$ws.Range("E26").Value = "'0.05%"
$ws.Range("G26").Value = "'19"
$ws.Range("G27").Value = "'19"
$ws.Range("G28").Value = "'19"
$ws.Range("G29").Value = "'19"
$ws.Range("G30").Value = "'19"
$ws.Range("G31").Value = "'19"
$ws.Range("G32").Value = "'19"
$ws.Range("G33").Value = "'19"
$ws.Range("G34").Value = "'19"
$ws.Range("G35").Value = "'19"
$ws.Range("G36").Value = "'19"
$ws.Range("G37").Value = "'19"
$ws.Range("D38").Value = "'0.02559"
$ws.Range("E38").Value = "'-6.71%"
$ws.Range("G38").Value = "'19"
$ws.Range("D39").Value = "'0.05419"
$ws.Range("E39").Value = "'-4.02%"
$ws.Range("G39").Value = "'19"
$ws.Range("D40").Value = "'0.008082"
$ws.Range("E40").Value = "'2.74%"
$ws.Range("G40").Value = "'19"
$ws.Range("D41").Value = "'0.1389"
$ws.Range("E41").Value = "'-1.71%"
$ws.Range("G41").Value = "'19"
$ws.Range("D42").Value = "'0.007581"
$ws.Range("E42").Value = "'2.82%"
$ws.Range("G42").Value = "'19"
$ws.Range("D43").Value = "'0.001987"
$ws.Range("E43").Value = "'-3.81%"
$ws.Range("G43").Value = "'19"
$ws.Range("D44").Value = "'0.008405"
$ws.Range("E44").Value = "'6.16%"
$ws.Range("G44").Value = "'19"
$ws.Range("D45").Value = "'0.3137"
$ws.Range("E45").Value = "'-10.27%"
$ws.Range("G45").Value = "'19"
$ws.Range("D46").Value = "'0.00006773"
$ws.Range("E46").Value = "'-2.08%"
$ws.Range("G46").Value = "'19"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("G47").Value = "'19"
$ws.Range("D48").Value = "'0.003397"
$ws.Range("E48").Value = "'-2.82%"
$ws.Range("G48").Value = "'19"
$ws.Range("E49").Value = "'16.20%"
$ws.Range("G49").Value = "'19"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("G50").Value = "'19"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.09%"
$ws.Range("G51").Value = "'19"
